$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the 17 new rows (1257-1273) of GPS tracking data as a 2D array.
$data = New-Object 'object[,]' 17,22
$data[0,0] = "Entrainement"
$data[0,1] = 46059
$data[0,2] = "Global"
$data[0,3] = "J-1"
$data[0,4] = "Malik Boussaid"
$data[0,5] = "right back"
$data[0,6] = "01:11:34"
$data[0,7] = 4.65
$data[0,8] = 0.18
$data[0,9] = 4.47
$data[0,10] = 0.13
$data[0,11] = 0.04
$data[0,12] = 0.03
$data[0,13] = 0
$data[0,14] = 4
$data[0,15] = 3.61
$data[0,16] = 27.77
$data[0,17] = 4.63
$data[0,18] = 33
$data[0,19] = 5
$data[0,20] = 22
$data[0,21] = 9
$data[1,0] = "Entrainement"
$data[1,1] = 46059
$data[1,2] = "Global"
$data[1,3] = "J-1"
$data[1,4] = "Mattheo Haon"
$data[1,5] = "right back"
$data[1,6] = "01:12:45"
$data[1,7] = 4.73
$data[1,8] = 0.15
$data[1,9] = 4.57
$data[1,10] = 0.13
$data[1,11] = 0.03
$data[1,12] = 0
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 3.86
$data[1,16] = 24.94
$data[1,17] = 4.76
$data[1,18] = 15
$data[1,19] = 3
$data[1,20] = 8
$data[1,21] = 1
$data[2,0] = "Entrainement"
$data[2,1] = 46059
$data[2,2] = "Global"
$data[2,3] = "J-1"
$data[2,4] = "Ilan Ihaddadene"
$data[2,5] = "center midfield"
$data[2,6] = "01:14:19"
$data[2,7] = 5.15
$data[2,8] = 0.08
$data[2,9] = 5.07
$data[2,10] = 0.08
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0
$data[2,14] = 0
$data[2,15] = 4.11
$data[2,16] = 20.44
$data[2,17] = 4.53
$data[2,18] = 13
$data[2,19] = 1
$data[2,20] = 9
$data[2,21] = 1
$data[3,0] = "N3 J15 VS OM (B)"
$data[3,1] = 46060
$data[3,2] = "Global"
$data[3,3] = "M"
$data[3,4] = "Ilan Ihaddadene"
$data[3,5] = "center midfield"
$data[3,6] = "00:12:38"
$data[3,7] = 1.67
$data[3,8] = 0.49
$data[3,9] = 1.17
$data[3,10] = 0.36
$data[3,11] = 0.12
$data[3,12] = 0.03
$data[3,13] = 0
$data[3,14] = 3
$data[3,15] = 7.88
$data[3,16] = 26.83
$data[3,17] = 4.84
$data[3,18] = 4
$data[3,19] = 1
$data[3,20] = 3
$data[3,21] = 1
$data[4,0] = "N3 J15 VS OM (B)"
$data[4,1] = 46060
$data[4,2] = "Global"
$data[4,3] = "M"
$data[4,4] = "Yoann Martelat"
$data[4,5] = "center midfield"
$data[4,6] = "01:27:00"
$data[4,7] = 10.84
$data[4,8] = 2.7
$data[4,9] = 8.11
$data[4,10] = 1.94
$data[4,11] = 0.68
$data[4,12] = 0.11
$data[4,13] = 0
$data[4,14] = 7
$data[4,15] = 7.44
$data[4,16] = 28.36
$data[4,17] = 4.65
$data[4,18] = 16
$data[4,19] = 2
$data[4,20] = 27
$data[4,21] = 6
$data[5,0] = "N3 J15 VS OM (B)"
$data[5,1] = 46060
$data[5,2] = "Global"
$data[5,3] = "M"
$data[5,4] = "Naim Dhib"
$data[5,5] = "center midfield"
$data[5,6] = "01:40:30"
$data[5,7] = 10.19
$data[5,8] = 1.7
$data[5,9] = 8.47
$data[5,10] = 1.34
$data[5,11] = 0.32
$data[5,12] = 0.06
$data[5,13] = 0
$data[5,14] = 4
$data[5,15] = 6.08
$data[5,16] = 28.03
$data[5,17] = 4.63
$data[5,18] = 42
$data[5,19] = 5
$data[5,20] = 45
$data[5,21] = 9
$data[6,0] = "N3 J15 VS OM (B)"
$data[6,1] = 46060
$data[6,2] = "Global"
$data[6,3] = "M"
$data[6,4] = "Kamal Bafounta"
$data[6,5] = "center midfield"
$data[6,6] = "00:26:53"
$data[6,7] = 3.36
$data[6,8] = 0.92
$data[6,9] = 2.43
$data[6,10] = 0.75
$data[6,11] = 0.17
$data[6,12] = 0.01
$data[6,13] = 0
$data[6,14] = 2
$data[6,15] = 7.43
$data[6,16] = 25.85
$data[6,17] = 5.09
$data[6,18] = 9
$data[6,19] = 4
$data[6,20] = 11
$data[6,21] = 5
$data[7,0] = "N3 J15 VS OM (B)"
$data[7,1] = 46060
$data[7,2] = "Global"
$data[7,3] = "M"
$data[7,4] = "Naim Ighbane"
$data[7,5] = "center back"
$data[7,6] = "01:40:23"
$data[7,7] = 10.21
$data[7,8] = 1.78
$data[7,9] = 8.4
$data[7,10] = 1.02
$data[7,11] = 0.6
$data[7,12] = 0.17
$data[7,13] = 0
$data[7,14] = 11
$data[7,15] = 6.05
$data[7,16] = 29.69
$data[7,17] = 4.16
$data[7,18] = 31
$data[7,19] = 4
$data[7,20] = 25
$data[7,21] = 11
$data[8,0] = "N3 J15 VS OM (B)"
$data[8,1] = 46060
$data[8,2] = "Global"
$data[8,3] = "M"
$data[8,4] = "Nathanael Beta"
$data[8,5] = "left forward"
$data[8,6] = "00:26:15"
$data[8,7] = 3.06
$data[8,8] = 0.57
$data[8,9] = 2.48
$data[8,10] = 0.38
$data[8,11] = 0.16
$data[8,12] = 0.05
$data[8,13] = 0
$data[8,14] = 6
$data[8,15] = 6.93
$data[8,16] = 27.65
$data[8,17] = 4.89
$data[8,18] = 15
$data[8,19] = 4
$data[8,20] = 14
$data[8,21] = 3
$data[9,0] = "N3 J15 VS OM (B)"
$data[9,1] = 46060
$data[9,2] = "Global"
$data[9,3] = "M"
$data[9,4] = "Yoan Zouma"
$data[9,5] = "center back"
$data[9,6] = "00:48:59"
$data[9,7] = 4.89
$data[9,8] = 0.68
$data[9,9] = 4.2
$data[9,10] = 0.49
$data[9,11] = 0.11
$data[9,12] = 0.08
$data[9,13] = 0.01
$data[9,14] = 5
$data[9,15] = 5.88
$data[9,16] = 31.31
$data[9,17] = 3.79
$data[9,18] = 9
$data[9,19] = 0
$data[9,20] = 12
$data[9,21] = 3
$data[10,0] = "N3 J15 VS OM (B)"
$data[10,1] = 46060
$data[10,2] = "Global"
$data[10,3] = "M"
$data[10,4] = "Mehdi Boussaid"
$data[10,5] = "center midfield"
$data[10,6] = "01:13:29"
$data[10,7] = 9.03
$data[10,8] = 2.43
$data[10,9] = 6.58
$data[10,10] = 1.59
$data[10,11] = 0.66
$data[10,12] = 0.2
$data[10,13] = 0
$data[10,14] = 14
$data[10,15] = 7.31
$data[10,16] = 28.75
$data[10,17] = 4.27
$data[10,18] = 33
$data[10,19] = 2
$data[10,20] = 34
$data[10,21] = 6
$data[11,0] = "N3 J15 VS OM (B)"
$data[11,1] = 46060
$data[11,2] = "Global"
$data[11,3] = "M"
$data[11,4] = "Karahali Souaré"
$data[11,5] = "right forward"
$data[11,6] = "00:48:59"
$data[11,7] = 5.79
$data[11,8] = 1.32
$data[11,9] = 4.45
$data[11,10] = 0.85
$data[11,11] = 0.41
$data[11,12] = 0.08
$data[11,13] = 0
$data[11,14] = 7
$data[11,15] = 7.04
$data[11,16] = 28.34
$data[11,17] = 5.11
$data[11,18] = 26
$data[11,19] = 7
$data[11,20] = 25
$data[11,21] = 16
$data[12,0] = "N3 J15 VS OM (B)"
$data[12,1] = 46060
$data[12,2] = "Global"
$data[12,3] = "M"
$data[12,4] = "Maé Clavel"
$data[12,5] = "left back"
$data[12,6] = "01:40:30"
$data[12,7] = 12.28
$data[12,8] = 2.48
$data[12,9] = 9.77
$data[12,10] = 1.62
$data[12,11] = 0.76
$data[12,12] = 0.13
$data[12,13] = 0
$data[12,14] = 10
$data[12,15] = 7.26
$data[12,16] = 28.69
$data[12,17] = 4.9
$data[12,18] = 47
$data[12,19] = 4
$data[12,20] = 45
$data[12,21] = 18
$data[13,0] = "N3 J15 VS OM (B)"
$data[13,1] = 46060
$data[13,2] = "Global"
$data[13,3] = "M"
$data[13,4] = "Sofiane Belle"
$data[13,5] = "left forward"
$data[13,6] = "01:39:44"
$data[13,7] = 10.25
$data[13,8] = 1.97
$data[13,9] = 8.25
$data[13,10] = 1.38
$data[13,11] = 0.52
$data[13,12] = 0.07
$data[13,13] = 0.03
$data[13,14] = 8
$data[13,15] = 6.11
$data[13,16] = 32.25
$data[13,17] = 4.4
$data[13,18] = 34
$data[13,19] = 4
$data[13,20] = 29
$data[13,21] = 12
$data[14,0] = "N3 J15 VS OM (B)"
$data[14,1] = 46060
$data[14,2] = "Global"
$data[14,3] = "M"
$data[14,4] = "Theo Owono"
$data[14,5] = "center midfield"
$data[14,6] = "00:44:11"
$data[14,7] = 5.11
$data[14,8] = 1.23
$data[14,9] = 3.86
$data[14,10] = 0.76
$data[14,11] = 0.38
$data[14,12] = 0.11
$data[14,13] = 0.01
$data[14,14] = 12
$data[14,15] = 6.86
$data[14,16] = 30.39
$data[14,17] = 4.95
$data[14,18] = 23
$data[14,19] = 6
$data[14,20] = 20
$data[14,21] = 18
$data[15,0] = "N3 J15 VS OM (B)"
$data[15,1] = 46060
$data[15,2] = "Global"
$data[15,3] = "M"
$data[15,4] = "Romain Thunet"
$data[15,5] = "center back"
$data[15,6] = "01:40:22"
$data[15,7] = 11.29
$data[15,8] = 2.1
$data[15,9] = 9.17
$data[15,10] = 1.43
$data[15,11] = 0.49
$data[15,12] = 0.19
$data[15,13] = 0.01
$data[15,14] = 18
$data[15,15] = 6.7
$data[15,16] = 30.41
$data[15,17] = 5.15
$data[15,18] = 50
$data[15,19] = 8
$data[15,20] = 37
$data[15,21] = 13
$data[16,0] = "N3 J15 VS OM (B)"
$data[16,1] = 46060
$data[16,2] = "Global"
$data[16,3] = "M"
$data[16,4] = "Jeremie Laurent"
$data[16,5] = "left forward"
$data[16,6] = "01:12:59"
$data[16,7] = 8.48
$data[16,8] = 2.24
$data[16,9] = 6.2
$data[16,10] = 1.47
$data[16,11] = 0.55
$data[16,12] = 0.24
$data[16,13] = 0.01
$data[16,14] = 11
$data[16,15] = 6.93
$data[16,16] = 30.64
$data[16,17] = 5.41
$data[16,18] = 54
$data[16,19] = 17
$data[16,20] = 36
$data[16,21] = 33

# Write the block of new data to A1257:V1273
$target = $ws.Range("A1257:V1273")
$target.Value = $data

# Column B holds dates; copy number format (style) from an existing date cell
# so the new cells render as dates (matches existing column B formatting).
$ws.Cells.Item(2, 2).Copy()
$ws.Range("B1257:B1273").PasteSpecial(-4122)

# Column D ("MD") uses a centered style; copy that format from an existing
# cell in the same column so the new entries match.
$ws.Cells.Item(1250, 4).Copy()
$ws.Range("D1257:D1273").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view: scroll position and active cell/selection, matching
# where the user was working after adding the new rows.
$excel.ActiveWindow.ScrollRow = 1245
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E1278").Select() | Out-Null
